$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.750.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "531.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.124.15"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.32"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.413"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.07%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.661.07"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.63%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "57.893.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.132.32"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.05"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "365.27"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.67"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.14"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.506"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0864"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.29"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.87"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.39"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.75%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "159.59"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.80%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.06"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.29"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.70%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.44"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.19%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0670"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.66%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.529.51"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.00%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.03"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.699"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "37.79"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0269"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.982"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.10"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.72"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.739"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.21%  "
